# Edit script: Cambios en el operador de Reparacion GACEP
# 1. Update AGEU450 row 2 values.
# 2. Add new sheet "GACEPReparacionMutacion" (copy of GACEPMutacion) with its own row 2 values.
# 3. Fix up selections / active sheet to match the target state.

$wb = $excel.ActiveWorkbook

$wsAGEU450 = $wb.Worksheets.Item("AGEU450")
$wsGACEP   = $wb.Worksheets.Item("GACEPMutacion")


# --- 1. Update AGEU450 (sheet1) row 2 values ---
$sheet1Row2 = @(192630,81790,188701,223666,228927,73290.5,10153,62355.5,232238,24746.4,18252.7,56447.6,3703.5,49946.1,61494,35953.8,14524.9,20280.6,35033.3,88004.1,83364.3,104608,33573.7,105713,56002.5,15949.6,52529.1,53674.4,68215.7,143452,188657,94671.5,61645.7,71876,27543.4,144869,110417,19381.6,104031,376115,936479,302088,29179.2,100548,779178,40468.6,693113,778717,622530,47941.4,202982,239154,220645,185475,78900.3,58001.3,147700,48885,281346,368682,209505,224676,225646,478885,425320,217170,315440,104041,142199,439128,283127,61507.5,126665,137097,227386,267436,596964,513040,374916,28698.3,275303,440074,14752.7,263446,481184,9116.6,244729,990479,500926,104813,866301,303098,715177,723710,43420.2,754010,752160)
for ($i = 0; $i -lt $sheet1Row2.Length; $i++) {
    $wsAGEU450.Cells.Item(2, $i + 1).Value = $sheet1Row2[$i]
}

# --- 2. Create GACEPReparacionMutacion as a copy of GACEPMutacion, placed right after it ---
$wsGACEP.Copy($null, $wsGACEP)
$wsNew = $wb.Worksheets.Item($wsGACEP.Index + 1)
$wsNew.Name = "GACEPReparacionMutacion"

# Drop the inherited column-A width override so the new sheet uses the default width.
$wsNew.Columns.Item(1).ClearFormats() | Out-Null

# --- 3. Write the new sheet's own row 2 values (overwrite what Copy() brought over) ---
$sheet3Row2 = @(212526,79531.7,213189,244658,248467,72115.1,13197.5,65598.8,248484,26370.2,20030.2,60477.8,4234.2,56899.4,64455.7,38102.4,14834.4,21250.8,37334.4,99631,87983.4,115638,36146.2,115302,55284.6,18582.2,54080,52731.7,72656.6,158357,192165,105867,61514.1,76063,29468.8,161562,121305,21361,119014,313430,983266,253149,32978.6,94315.9,825296,44124.1,717162,816448,632838,39488.7,206967,249245,233681,191942,66763.3,50194.9,137753,39433.7,273776,367681,187042,191384,196566,499320,448356,197810,311635,72945.5,116643,431662,245982,53478.3,95893,109840,195521,229155,608935,512333,365222,24338.7,242051,439120,14305.4,225621,491344,9796.2,198785,1033690,426581,67840.6,885358,221735,694994,698756,41561.2,729626,725852)
for ($i = 0; $i -lt $sheet3Row2.Length; $i++) {
    $wsNew.Cells.Item(2, $i + 1).Value = $sheet3Row2[$i]
}

# --- 4. Restore the expected selections on each sheet ---
$wsAGEU450.Range("A2").Select() | Out-Null
$wsGACEP.Range("A1").Select() | Out-Null
$wsNew.Range("A4").Select() | Out-Null

# --- 5. Make the new sheet the active tab (tabSelected / activeTab) ---
$wsNew.Activate() | Out-Null

Write-Host "GACEPReparacionMutacion created; sheets now:"
foreach ($s in $wb.Worksheets) { Write-Host " -" $s.Name }
